$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1412
$unit = 48

$data = @(
    @("U48_01", 48, "Xuất bản", "Publish", "The author published one novel", "publish something ", "V"),
    @("U48_02", 48, "Người máy", "Robot", "I always sing without emotion like a robot", "like a robot / như một cái máy", "N"),
    @("U48_03", 48, "Máy móc", "Machine", "The watch was made by machine", "by machine / bằng máy", "N"),
    @("U48_04", 48, "Có thể", "able", "She is able to play the piano", "able to do something / có thể làm một cái gì đó", "Adj"),
    @("U48_05", 48, "Tự động", "Automatic", "The car is fully automatic", "fully automatic / hoàn toàn tự động", "Adj"),
    @("U48_06", 48, "Trả lời", "Respond", "Please respond to my email by morning", "respond to something or somebody ", "V"),
    @("U48_07", 48, "Khảo sát", "Survey", "We conducted a survey for research", "conduct a survey / tiến hành một cuộc khảo sát", "N"),
    @("U48_08", 48, "Dòng trạng thái", "Status", "I recently posted a status on FB", "post a status / đăng tải một dòng trạng thái", "N"),
    @("U48_09", 48, "Đính", "Attach", "Could you attach the file with your email?", "attach something to something / đính một cái gì đó vào một cái gì đó", "V"),
    @("U48_10", 48, "Tài khoản", "Account", "I am opening up a bank account", "a bank account", "N"),
    @("U48_11", 48, "Chính sách", "Policy", "What is your policy on education?", "policy on something / chính sách về một cái gì đó", "N"),
    @("U48_12", 48, "Trang", "Page", "Turn to page 25 and read", "turn to page / lật sang trang", "N"),
    @("U48_13", 48, "Cho phép", "Permit", "My parents do not permit me to come home late", "permit somebody to do something / cho pe ai đó làm một cái gì đó", "V"),
    @("U48_14", 48, "Hợp lệ, có hiệu lực", "Valid", "A valid passport is needed for traveling", "valid passport / hộ chiếu có hiệu lực", "Adj"),
    @("U48_15", 48, "Ổn định", "Stable", "The petrol price remains stable in recent times", "remain stable / vẫn ổn định", "Adj"),
    @("U48_16", 48, "Kẻ trộm", "Thief", "Stealing makes you a petty thief", "petty thief / kẻ cắp vặt", "N"),
    @("U48_17", 48, "Vệ tinh", "Satellite", "We can watch TV via satallite", "via satellite / thông qua vệ tinh", "N"),
    @("U48_18", 48, "Tên lửa", "Rocket", "Little boys love space rockets", "space rocket / tên lửa vũ trụ", "N"),
    @("U48_19", 48, "Vũ khí", "Weapon", "Police are armed (trang bị) with weapons", "with a weapon / với một vũ khí", "N"),
    @("U48_20", 48, "Lời cảnh cáo, cảnh báo", "Warning", "I received a warning for my mistake", "receive a warning / nhận được một lời cảnh cáo", "N"),
    @("U48_21", 48, "Hồ sơ", "Profile", "I set up a profile on FB", "set up a profile / thiết lập hồ sơ cá nhân", "N"),
    @("U48_22", 48, "Chọn", "Select", "I will help you select an option for your travel plan", "select something / lựa chọn một cái gì đó", "V"),
    @("U48_23", 48, "Chỉnh sửa", "Edit", "Can you edit my essay?", "edit something / chỉnh sửa một cái gì đó", "V"),
    @("U48_24", 48, "Giới hạn", "Limit", "Mother limited me to playing once a day.", "limit somebody to something / giới hạn ai đó trong một việc gì đó", "V"),
    @("U48_25", 48, "Sao lưu dữ liệu", "Back up", "You should back up your data", "back up the data", "V"),
    @("U48_26", 48, "Dữ liệu", "Data", "Scientists gather data for research ", "gather data / thu thập dữ liệu", "N"),
    @("U48_27", 48, "Thư viện", "Library", "Return the library book in 1 day", "library book / sách thư viện", "N"),
    @("U48_28", 48, "Tìm kiếm", "Search", "I search for love", "search for somebody or something", "V"),
    @("U48_29", 48, "Tiêu chuẩn", "Standard", "All products must meet the quality standards", "quality standards / các tiêu chuẩn chất lượng", "N"),
    @("U48_30", 48, "Kho", "Storage", "The bags are still in storage", "in storage / trong kho", "N")
)

# Pass 1: write all the new IDs (column A) first, matching the order
# in which the original authoring tool allocated shared-string entries.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
}

# Pass 2: write the remaining columns row by row.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$win = $excel.ActiveWindow
$win.ScrollRow = 1436
$win.Zoom = 120

$ws.Range("C1455").Select()
